$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text type (avoid Excel auto-converting
# numeric-looking strings like "584.30" or "0.540" into Number values,
# which would silently drop significant trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '74.622.77'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +8.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.595.34'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.87%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +15.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '584.30'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.540'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.207'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +23.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.593.20'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.87%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.363'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +10.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.81'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.60%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +9.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.525.11'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +8.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.070.93'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.89%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +13.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.610.69'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +7.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.25'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +33.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.82'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +12.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +12.67%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +18.64%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.48%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.55'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.19'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +14.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.33'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.739.47'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0952'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +16.68%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +20.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.97'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +12.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '508.95'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +19.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.76'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +9.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.121'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +15.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.97'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.26'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.35'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.69%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.93'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.65%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.69'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +12.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.328'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +10.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.43'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +19.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '156.54'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +18.66%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.79'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0847'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +18.09%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.80%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +9.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.29'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +20.38%  '
